# Reorders the "Recorded By" column (G) so that the literal entry "System"
# (exact case) always appears first in the comma-separated list of
# recorder names, while preserving the relative order of the remaining
# entries. Rows where "System" is absent, or already first, are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $raw = $cell.Text

    if ([string]::IsNullOrEmpty($raw)) {
        continue
    }

    $parts = $raw -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $systemIndex = [array]::IndexOf($parts, 'System')

    if ($systemIndex -gt 0) {
        $reordered = New-Object System.Collections.Generic.List[string]
        [void]$reordered.Add('System')
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $systemIndex) {
                [void]$reordered.Add($parts[$i])
            }
        }
        $newValue = [string]::Join(', ', $reordered)
        $cell.Value2 = $newValue
    }
}
